$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 11:35"

# Update province A22/A23 labels (La Rioja now has more cases than Salamanca, so they swap rank)
$ws.Range("A22").Value = "La Rioja"
$ws.Range("A23").Value = "Salamanca"

# --- Update numeric data (B:E) for the changed rows ---

# Row 4 - Madrid
$ws.Range("B4").Value = 65693
$ws.Range("C4").Value = 40199
$ws.Range("D4").Value = 16715
$ws.Range("E4").Value = 8779

# Row 5 - Cataluña
$ws.Range("B5").Value = 55482
$ws.Range("C5").Value = 25849
$ws.Range("D5").Value = 23810
$ws.Range("E5").Value = 5823

# Row 6 - Castilla y Leon
$ws.Range("B6").Value = 18173
$ws.Range("C6").Value = 7621
$ws.Range("D6").Value = 8627
$ws.Range("E6").Value = 1925

# Row 7 - Castilla-La Mancha
$ws.Range("B7").Value = 16470
$ws.Range("C7").Value = 6244
$ws.Range("D7").Value = 7374
$ws.Range("E7").Value = 2852

# Row 9 - Andalucia
$ws.Range("B9").Value = 12359
$ws.Range("C9").Value = 9918
$ws.Range("D9").Value = 1105
$ws.Range("E9").Value = 1336

# Row 11 - Galicia
$ws.Range("B11").Value = 9317
$ws.Range("C11").Value = 8157
$ws.Range("D11").Value = 559
$ws.Range("E11").Value = 601

# Row 14 - Aragon
$ws.Range("B14").Value = 5389
$ws.Range("C14").Value = 3471
$ws.Range("D14").Value = 1082
$ws.Range("E14").Value = 836

# Row 16 - Navarra
$ws.Range("B16").Value = 5105
$ws.Range("C16").Value = 3400
$ws.Range("D16").Value = 1207
$ws.Range("E16").Value = 498

# Row 22 - now La Rioja
$ws.Range("B22").Value = 4014
$ws.Range("C22").Value = 2867
$ws.Range("D22").Value = 799
$ws.Range("E22").Value = 348

# Row 23 - now Salamanca
$ws.Range("B23").Value = 4012
$ws.Range("C23").Value = 1124
$ws.Range("D23").Value = 2541
$ws.Range("E23").Value = 347

# Row 33 - Asturias (B33 unchanged)
$ws.Range("C33").Value = 1046
$ws.Range("D33").Value = 1002
$ws.Range("E33").Value = 308

# Row 34 - Gran Canaria (E34 unchanged)
$ws.Range("B34").Value = 2275
$ws.Range("C34").Value = 1496
$ws.Range("D34").Value = 628
